$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Regulation value was re-entered in lowercase ("csvla" instead of "CSVLA")
$ws.Range("B5").Value = "csvla"

# Update the data-validation list to match the new lowercase wording
$ws.Range("B5").Validation.Modify(3, 1, 1, '"csvla, cs-23, cs-25"')

# Reflect the scrolled/selected view state recorded in the saved file
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B31").Select()
